$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.699.80"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.894.03"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "242.09"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "0.06744"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("D10").Value = "1.894.39"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "17.17"
$ws.Range("E11").Value = "  +5.52%  "
$ws.Range("D12").Value = "0.07262"
$ws.Range("D13").Value = "90.86"
$ws.Range("E13").Value = "  +5.78%  "
$ws.Range("D14").Value = "0.6759"
$ws.Range("E14").Value = "  +2.01%  "
$ws.Range("D15").Value = "5.032"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "30.671.09"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "0.000007972"
$ws.Range("E17").Value = "  +2.73%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "2.140.68"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "4.802"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "191.28"
$ws.Range("E23").Value = "  +33.50%  "
$ws.Range("D24").Value = "6.092"
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("D25").Value = "9.371"
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("D26").Value = "157.47"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  +11.37%  "
$ws.Range("D28").Value = "1.897"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").Value = "4.303"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").Value = "0.09093"
$ws.Range("D32").Value = "4.005"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "0.05239"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "0.7396"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").Value = "1.105"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "2.738"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").Value = "0.01828"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").Value = "0.9352"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "2.126"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "0.4405"
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("D42").Value = "105.09"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.1351"
$ws.Range("E45").Value = "  +5.29%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.522"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("D47").Value = "0.05871"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "8.750"
$ws.Range("E48").Value = "  +6.20%  "
$ws.Range("D49").Value = "1.428"
$ws.Range("E49").Value = "  +6.25%  "
$ws.Range("D50").Value = "33.81"
$ws.Range("D51").Value = "0.3943"
$ws.Range("E51").Value = "  +5.00%  "
